$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-5 from 45243 to 45244
# (2023-11-13 -> 2023-11-14), keeping existing date formatting/style.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45244
}
